# Add "I0" (column I) and "IF" (column J) next to the existing "IP" (column H)
# column. I0 is a constant multiplier of 1; IF mirrors the IP value (IP * I0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting from the existing "IP" header cell (H1) so
# the new headers match the rest of the header row (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-39: I is always 1, J mirrors the H (IP) value for that row.
for ($row = 2; $row -le 39; $row++) {
    $ipValue = $ws.Cells.Item($row, 8).Value2

    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $ipValue
}
